$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 2 (Target cluster -> FAPs)
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 1.315741666666667
$ws.Range("H2").Value = 3.947225
$ws.Range("M2").Value = 198.2465873333333
$ws.Range("N2").Value = 594.7397619999999
$ws.Range("O2").Value = 0.9851515664921635
$ws.Range("P2").Value = 0.9851515664921635
$ws.Range("Q2").Value = 260.8412952289389
$ws.Range("R2").Value = 2347.57165706045
$ws.Range("S2").Value = 0.9851515664921635
$ws.Range("T2").Value = 0.9851515664921635

# New values for row 3 (Target cluster -> MuSCs)
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 1.315741666666667
$ws.Range("H3").Value = 3.947225
$ws.Range("M3").Value = 2.988018666666667
$ws.Range("N3").Value = 8.964055999999999
$ws.Range("O3").Value = 0.01484843350783645
$ws.Range("P3").Value = 0.01484843350783645
$ws.Range("Q3").Value = 3.931460660511111
$ws.Range("R3").Value = 35.3831459446
$ws.Range("S3").Value = 0.01484843350783645
$ws.Range("T3").Value = 0.01484843350783645

# Delete the old row 4 entirely (data moved into rows 2/3 above)
$ws.Rows(4).Delete()
